$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.927.18"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.523.34"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "587.89"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "171.46"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "2.523.42"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "2.995.63"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").Value = "66.837.23"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "2.531.81"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +3.04%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "353.15"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  -1.64%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +4.66%  "
$ws.Range("E25").Value = "  +0.00%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "69.72"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "2.682.44"
$ws.Range("E28").Value = "  -0.95%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.989"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -1.23%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "531.87"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.72%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.08"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("E35").Value = "  -0.81%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.45"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "156.37"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.34%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.54"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("E45").Value = "  +3.25%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "39.73"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "148.66"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "0.0₆0276"
$ws.Range("E49").Value = "  -4.11%  "
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("E51").Value = "  -0.04%  "
